$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the "Effort" (column F) values for the existing backlog rows ---
$ws.Range("F2").Value = 8
$ws.Range("F3").Value = 8
$ws.Range("F4").Value = 8
$ws.Range("F5").Value = 13
$ws.Range("F6").Value = 8
$ws.Range("F7").Value = 20
$ws.Range("F8").Value = 40
$ws.Range("F9").Value = 13
$ws.Range("F10").Value = 20
$ws.Range("F11").Value = 8
$ws.Range("F12").Value = 20
$ws.Range("F13").Value = 5
$ws.Range("F14").Value = 20
$ws.Range("F15").Value = 13
$ws.Range("F16").Value = 8
$ws.Range("F17").Value = 3
$ws.Range("F18").Value = 13

# --- New backlog items (rows 19-21) ---
$ws.Range("B19").Value = "Implement speed slider"
$ws.Range("C19").Value = "Low"
$ws.Range("D19").Value = "Sprint ready"
$ws.Range("E19").Value = "Implement slider to control bpm of the song"
$ws.Range("F19").Value = 13

$ws.Range("B20").Value = "Problems with wav file format"
$ws.Range("C20").Value = "Medium"
$ws.Range("D20").Value = "Sprint ready"
$ws.Range("E20").Value = "Should implement normalizer of wav files"
$ws.Range("F20").Value = 13

$ws.Range("B21").Value = "Too much clicking on track to add wav file cause System.ArgumentOutOfRangeException"
$ws.Range("C21").Value = "Medium"
$ws.Range("D21").Value = "Sprint ready"
$ws.Range("E21").Value = "Probably because dobule clicked on same position"
$ws.Range("F21").Value = 8

# --- Row height for the wrapped / taller row 21 ---
$ws.Rows.Item(21).RowHeight = 30

# --- Column width adjustments ---
$ws.Columns.Item(2).ColumnWidth = 56.3
$ws.Columns.Item(6).ColumnWidth = 5.14

# --- Selection / view state ---
$ws.Range("F21").Select()
